$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 0.3694329494838015
$ws.Range("F2").Value = 0.9983038269360573
$ws.Range("G2").Value = -0.4652015470543536
$ws.Range("H2").Value = -0.6458914122405494
$ws.Range("I2").Value = -47.80761937369478
$ws.Range("J2").Value = 0.5573518019022105
$ws.Range("K2").Value = 0.5225679220907508
$ws.Range("L2").Value = 0.7641188735158745
$ws.Range("M2").Value = 0.6422017199765734
$ws.Range("B3").Value = 4
$ws.Range("E3").Value = 0.3579305072448689
$ws.Range("F3").Value = 0.8810055940559142
$ws.Range("G3").Value = -0.3085431251618616
$ws.Range("H3").Value = -0.6833757971094011
$ws.Range("I3").Value = -43.14622557403036
$ws.Range("J3").Value = 0.5410777066202382
$ws.Range("K3").Value = 0.4916271061158151
$ws.Range("L3").Value = 0.7310699129959882
$ws.Range("M3").Value = 0.4893146111806932
$ws.Range("B4").Value = 7
$ws.Range("E4").Value = 0.3440443284603508
$ws.Range("F4").Value = 0.8600973034666924
$ws.Range("G4").Value = -0.2329571306089141
$ws.Range("H4").Value = -0.6400746383872615
$ws.Range("I4").Value = -37.68503496879943
$ws.Range("J4").Value = 0.5341894422128338
$ws.Range("K4").Value = 0.4997862701244386
$ws.Range("L4").Value = 0.7315433308844582
$ws.Range("M4").Value = 0.4521664697450679
$ws.Range("B5").Value = 10
$ws.Range("E5").Value = 0.3134214201635895
$ws.Range("F5").Value = 0.7923424090903093
$ws.Range("G5").Value = -0.1884625918330206
$ws.Range("H5").Value = -0.5783821867279736
$ws.Range("I5").Value = -32.98401933015647
$ws.Range("J5").Value = 0.5447701960051228
$ws.Range("K5").Value = 0.5006896374233462
$ws.Range("L5").Value = 0.7399284091330687
$ws.Range("M5").Value = 0.283077820693207
$ws.Range("B6").Value = 13
$ws.Range("E6").Value = 0.2834625187009352
$ws.Range("F6").Value = 0.7717089551300679
$ws.Range("G6").Value = -0.05053952486684921
$ws.Range("H6").Value = -0.5524418853559041
$ws.Range("I6").Value = -28.4098335242593
$ws.Range("J6").Value = 0.5390622867544849
$ws.Range("K6").Value = 0.4991468719885641
$ws.Range("L6").Value = 0.7346702024092181
$ws.Range("M6").Value = 0.2359148008558608
$ws.Range("B7").Value = 16
$ws.Range("E7").Value = 0.2951024894186304
$ws.Range("F7").Value = 0.7738981661897398
$ws.Range("G7").Value = -0.02974457688412713
$ws.Range("H7").Value = -0.391565488310808
$ws.Range("I7").Value = -22.86775401229214
$ws.Range("J7").Value = 0.5361164356316949
$ws.Range("K7").Value = 0.4930037230054645
$ws.Range("L7").Value = 0.728339205312858
$ws.Range("M7").Value = 0.2614213927568656
$ws.Range("B8").Value = 19
$ws.Range("E8").Value = 0.261878499449297
$ws.Range("F8").Value = 0.7382036408632925
$ws.Range("G8").Value = 0.03829161117804839
$ws.Range("H8").Value = -0.2444748509359442
$ws.Range("I8").Value = -18.20837363922245
$ws.Range("J8").Value = 0.4976729420960331
$ws.Range("K8").Value = 0.4662748812435507
$ws.Range("L8").Value = 0.6819872368735078
$ws.Range("M8").Value = 0.2824908349574051
$ws.Range("B9").Value = 22
$ws.Range("E9").Value = 0.2647974373044657
$ws.Range("F9").Value = 0.7497131728124901
$ws.Range("G9").Value = 0.009887648395874749
$ws.Range("H9").Value = -0.2404089343539037
$ws.Range("I9").Value = -13.52116639618493
$ws.Range("J9").Value = 0.5030890420873
$ws.Range("K9").Value = 0.4439084756232666
$ws.Range("L9").Value = 0.6710036004685098
$ws.Range("M9").Value = 0.3341665746314018
$ws.Range("B10").Value = 25
$ws.Range("E10").Value = 0.2032466775022916
$ws.Range("F10").Value = 0.6545653607509401
$ws.Range("G10").Value = 0.06974547826215864
$ws.Range("H10").Value = -0.1387982665010365
$ws.Range("I10").Value = -8.4048673398409
$ws.Range("J10").Value = 0.4795351001829849
$ws.Range("K10").Value = 0.3880976376770067
$ws.Range("L10").Value = 0.616921911796776
$ws.Range("M10").Value = 0.2187239530096449
$ws.Range("B11").Value = 28
$ws.Range("E11").Value = 0.240764652535601
$ws.Range("F11").Value = 0.6408152324475869
$ws.Range("G11").Value = -0.03524418780312999
$ws.Range("H11").Value = -0.05779066989024765
$ws.Range("I11").Value = -4.162544364945773
$ws.Range("J11").Value = 0.4450505049146024
$ws.Range("K11").Value = 0.3273556607150989
$ws.Range("L11").Value = 0.5524780150041632
$ws.Range("M11").Value = 0.3245067111361079
$ws.Range("B12").Value = 31
$ws.Range("E12").Value = 0.8538574628745913
$ws.Range("F12").Value = 1.389595833084475
$ws.Range("G12").Value = 0.01807536981090152
$ws.Range("H12").Value = -0.0006114703069025987
$ws.Range("I12").Value = 2.792443443671765
$ws.Range("J12").Value = 0.4974175322170247
$ws.Range("K12").Value = 0.437606130464236
$ws.Range("L12").Value = 0.6625130559929254
$ws.Range("M12").Value = 1.22149154391113
$ws.Range("B13").Value = 34
$ws.Range("E13").Value = 0.2437425433792841
$ws.Range("F13").Value = 0.6873233565783542
$ws.Range("G13").Value = -0.01025735848063884
$ws.Range("H13").Value = -0.01542737967620648
$ws.Range("I13").Value = 6.882243475509802
$ws.Range("J13").Value = 0.4940560442916189
$ws.Range("K13").Value = 0.4462672081408541
$ws.Range("L13").Value = 0.665768082280209
$ws.Range("M13").Value = 0.1706821646377162
$ws.Range("B14").Value = 37
$ws.Range("E14").Value = 0.2644352771148299
$ws.Range("F14").Value = 0.7164393154644291
$ws.Range("G14").Value = -0.01490295508841655
$ws.Range("I14").Value = 12.29418718127211
$ws.Range("J14").Value = 0.5213689386940633
$ws.Range("K14").Value = 0.4744964756449879
$ws.Range("L14").Value = 0.7049628013701005
$ws.Range("M14").Value = 0.127708006090825
$ws.Range("B15").Value = 40
$ws.Range("E15").Value = 0.2721814889639132
$ws.Range("F15").Value = 0.7246924555353139
$ws.Range("G15").Value = -0.04943590905526207
$ws.Range("H15").Value = 0.04113079035132235
$ws.Range("I15").Value = 17.00710186105163
$ws.Range("J15").Value = 0.5122924918626005
$ws.Range("K15").Value = 0.4962966200674679
$ws.Range("L15").Value = 0.7132700866563958
$ws.Range("M15").Value = 0.1280679253679703
$ws.Range("B16").Value = 43
$ws.Range("E16").Value = 0.2748700687140564
$ws.Range("F16").Value = 0.7406040875179626
$ws.Range("G16").Value = -0.1080439519744042
$ws.Range("H16").Value = 0.02435961947621005
$ws.Range("I16").Value = 22.14348471151862
$ws.Range("J16").Value = 0.5275406360953262
$ws.Range("K16").Value = 0.4954467018880422
$ws.Range("L16").Value = 0.7237194608814925
$ws.Range("M16").Value = 0.1571629583812731
$ws.Range("B17").Value = 46
$ws.Range("E17").Value = 0.2826214883730456
$ws.Range("F17").Value = 0.7539212910134512
$ws.Range("G17").Value = -0.1686652107780825
$ws.Range("H17").Value = 0.04348353214743383
$ws.Range("I17").Value = 27.84532288491843
$ws.Range("J17").Value = 0.5411950079529787
$ws.Range("K17").Value = 0.4916870139963372
$ws.Range("L17").Value = 0.7312181470797602
$ws.Range("M17").Value = 0.1835875192217451
$ws.Range("B18").Value = 49
$ws.Range("E18").Value = 0.2857128696952811
$ws.Range("F18").Value = 0.7581659534554785
$ws.Range("G18").Value = -0.2128307441399026
$ws.Range("H18").Value = 0.03217625927849591
$ws.Range("I18").Value = 32.63276593349078
$ws.Range("J18").Value = 0.5410547603474684
$ws.Range("K18").Value = 0.4869007103688891
$ws.Range("L18").Value = 0.7278962879208942
$ws.Range("M18").Value = 0.2115099193145234
$ws.Range("B19").Value = 52
$ws.Range("E19").Value = 0.2496752549320211
$ws.Range("F19").Value = 0.7555538448824658
$ws.Range("G19").Value = -0.2741783046714848
$ws.Range("H19").Value = 0.07084085213352864
$ws.Range("I19").Value = 37.69309004095294
$ws.Range("J19").Value = 0.5297601032149251
$ws.Range("K19").Value = 0.461872761248127
$ws.Range("L19").Value = 0.7029211875634135
$ws.Range("M19").Value = 0.2767795962419354
$ws.Range("B20").Value = 55
$ws.Range("E20").Value = 0.2820956853628623
$ws.Range("F20").Value = 0.7726528415011927
$ws.Range("G20").Value = -0.20520759521321
$ws.Range("H20").Value = 0.1993124582297507
$ws.Range("I20").Value = 42.43904040485148
$ws.Range("J20").Value = 0.5285290055298583
$ws.Range("K20").Value = 0.472053754396442
$ws.Range("L20").Value = 0.7087189697665247
$ws.Range("M20").Value = 0.3066845007229721
$ws.Range("B21").Value = 58
$ws.Range("E21").Value = 0.3300797486388904
$ws.Range("F21").Value = 0.8360733178083956
$ws.Range("G21").Value = -0.2853772018690013
$ws.Range("H21").Value = 0.110147094428991
$ws.Range("I21").Value = 47.73737586310463
$ws.Range("J21").Value = 0.5182003417486727
$ws.Range("K21").Value = 0.4722740768951236
$ws.Range("L21").Value = 0.7014998210886193
$ws.Range("M21").Value = 0.4537722706699222
$ws.Range("B22").Value = 61
$ws.Range("E22").Value = 0.4175725031442469
$ws.Range("F22").Value = 0.9932962400361531
$ws.Range("G22").Value = -0.3110882019058749
$ws.Range("H22").Value = 0.1971740793524835
$ws.Range("I22").Value = 53.27193727877913
$ws.Range("J22").Value = 0.5423428319286113
$ws.Range("K22").Value = 0.5247630533923773
$ws.Range("L22").Value = 0.754821414638457
$ws.Range("M22").Value = 0.6453518759482443
